$d = $word.ActiveDocument

$p69 = $d.Paragraphs.Item(69)
# sanity check anchor text
if ($p69.Range.Text -notmatch "\x07?$") { }

# (A1) After the bookmark paragraph: insert empty, CONCLUSION, empty
$p69.Range.InsertParagraphAfter() | Out-Null
$ph = $p69.Next()
$phRange = $ph.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t>CONCLUSION</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$phRange.InsertXML($xml)

# (A2) Prepend two new runs of text before the bookmarks inside paragraph 69
$p69 = $d.Paragraphs.Item(69)
$collapsedStart = $d.Range($p69.Range.Start, $p69.Range.Start)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">However, this implementation performs quite well, despite this lack of independence.  This might be due to the fact that we’re looking at somewhat distinct digits; all of the 7s look similar and not much like other digits.  Further, this is a bitmap; a pixel is either filled or not.  Therefore, </w:t></w:r><w:r><w:t>there isn’t any “grey area” (literally) to introduce error on fringe pixels.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$collapsedStart.InsertXML($xml)

# (B1) Replace paragraph 68's (old CONCLUSION) content with 3 new runs
$p68 = $d.Paragraphs.Item(68)
$p68RangeNoMark = $d.Range($p68.Range.Start, $p68.Range.End - 1)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Basic naïve Bayes performed rather well on this data set, scoring </w:t></w:r><w:r><w:t>almost 98% accuracy.</w:t></w:r><w:r><w:t xml:space="preserve">  I do not believe that these features are necessarily completely independent of each other.  The probability that a given pixel is filled would increase with the number of adjacent pixels also filled.  For example, when observing a 3x3 pixel matrix, if all edge pixels are filled, it would be very highly likely the center pixel is also filled.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p68RangeNoMark.InsertXML($xml)

# (B2) After paragraph 68: insert one new empty paragraph
$p68.Range.InsertParagraphAfter() | Out-Null
$ph = $p68.Next()
$phRange = $ph.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$phRange.InsertXML($xml)

# (C) After RESULTS (paragraph 64): insert the 17-paragraph results block
$d.Paragraphs.Item(64).Range.InsertParagraphAfter() | Out-Null
$ph = $d.Paragraphs.Item(64).Next()
$phRange = $ph.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t xml:space="preserve">Naïve Bayes </w:t></w:r><w:r><w:t xml:space="preserve">Test </w:t></w:r><w:r><w:t>D</w:t></w:r><w:r><w:t>ata:</w:t></w:r></w:p><w:p><w:r><w:t>1797 total instances</w:t></w:r></w:p><w:p><w:r><w:t>Total Accuracy: 97.94%</w:t></w:r></w:p><w:p><w:r><w:t>0: 173(TP), 1617(TN), 2(FP), 5(FN)</w:t></w:r></w:p><w:p><w:r><w:t>1: 153(TP), 1585(TN), 30(FP), 29(FN)</w:t></w:r></w:p><w:p><w:r><w:t>2: 152(TP), 1601(TN), 19(FP), 25(FN)</w:t></w:r></w:p><w:p><w:r><w:t>3: 158(TP), 1600(TN), 14(FP), 25(FN)</w:t></w:r></w:p><w:p><w:r><w:t>4: 170(TP), 1595(TN), 21(FP), 11(FN)</w:t></w:r></w:p><w:p><w:r><w:t>5: 166(TP), 1602(TN), 13(FP), 16(FN)</w:t></w:r></w:p><w:p><w:r><w:t>6: 176(TP), 1611(TN), 5(FP), 5(FN)</w:t></w:r></w:p><w:p><w:r><w:t>7: 169(TP), 1601(TN), 17(FP), 10(FN)</w:t></w:r></w:p><w:p><w:r><w:t>8: 139(TP), 1605(TN), 18(FP), 35(FN)</w:t></w:r></w:p><w:p><w:r><w:t>9: 156(TP), 1571(TN), 46(FP), 24(FN)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Binning Test Data:</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$phRange.InsertXML($xml)

# (D) After the blank line following OVERVIEW (paragraph 63): insert the 8-paragraph overview/implementation block
$d.Paragraphs.Item(63).Range.InsertParagraphAfter() | Out-Null
$ph = $d.Paragraphs.Item(63).Next()
$phRange = $ph.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>The computing portion of this project explored naïve Bayes classification of numerical digits, represented by 64 features indicating a bitmap of black and white pixels from each digit’s image.  Further, a binning technique was used to test for improved accuracy.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>IMPLEMENTATION</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">This project was coded in Java </w:t></w:r><w:r><w:t>1.7.0_09</w:t></w:r><w:r><w:t>, using Eclipse on a MacBook Pro running Mavericks.  Training and test data was stored in external files, read in and assimilated into respective arrays.  Probabilities were calculated for each digit, based on the number of occurrences within the training data, and for each digit’s feature, smoothed with Laplace smoothing.  Calculating all of these probabilities and overlaying these with the actual features of each test instance produced a probability for classification.  The digit classification with the highest probability was then selected as the test instance’s classification.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>For binning, the 17 values for each feature were reduced to four bins, abstracting the data and simplifying overall calculations.  This was implemented by simply reassigning feature values to a bin from 0-3, based on the original feature value.  Classification then proceeded as described above.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$phRange.InsertXML($xml)

# (E1) Add <w:lastRenderedPageBreak/> before OVERVIEW's text
$p62 = $d.Paragraphs.Item(62)
$p62RangeNoMark = $d.Range($p62.Range.Start, $p62.Range.End - 1)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>OVERVIEW</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p62RangeNoMark.InsertXML($xml)

# (E2) Insert a new paragraph containing a page break, immediately before OVERVIEW
$p61 = $d.Paragraphs.Item(61)
$p61.Range.InsertParagraphAfter() | Out-Null
$ph2 = $p61.Next()
$ph2Range = $ph2.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:br w:type="page"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$ph2Range.InsertXML($xml)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
